$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Feb 11 20:19:34 EST 2025"
$ws.Range("B3").Value = "Tue Feb 11 20:19:47 EST 2025"
$ws.Range("B4").Value = "Tue Feb 11 20:20:00 EST 2025"
$ws.Range("B5").Value = "Tue Feb 11 20:20:13 EST 2025"
$ws.Range("B6").Value = "Tue Feb 11 20:20:26 EST 2025"
$ws.Range("B7").Value = "Tue Feb 11 20:20:39 EST 2025"
